$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-28 Saturday" "2025-06-29 Sunday"

Replace-Text "30×44=1320" "85×46=3910"
Replace-Text "46×24=1104" "59×98=5782"
Replace-Text "35×92=3220" "42×48=2016"
Replace-Text "96×17=1632" "17×57=969"
Replace-Text "29×27=783" "11×65=715"

Replace-Text "86×51=4386" "93×15=1395"
Replace-Text "42×61=2562" "49×55=2695"
Replace-Text "35×70=2450" "60×60=3600"
Replace-Text "15×82=1230" "80×26=2080"
Replace-Text "27×64=1728" "24×74=1776"

Replace-Text "38×17=646" "99×33=3267"
Replace-Text "69×33=2277" "48×76=3648"
Replace-Text "58×81=4698" "41×74=3034"
Replace-Text "96×34=3264" "18×47=846"
Replace-Text "78×28=2184" "34×86=2924"

Replace-Text "55×50=2750" "90×91=8190"
Replace-Text "96×67=6432" "53×22=1166"
Replace-Text "23×87=2001" "76×81=6156"
Replace-Text "47×11=517" "48×33=1584"
Replace-Text "20×59=1180" "68×85=5780"

Replace-Text "20×28=560" "99×93=9207"
Replace-Text "34×44=1496" "38×95=3610"
Replace-Text "87×45=3915" "58×42=2436"
Replace-Text "68×80=5440" "63×67=4221"
Replace-Text "56×90=5040" "66×99=6534"
